$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = "my_files\קורות חיים מיכל לאער (4).pdf"
$ws.Range("B26").Value = "בתקשורת -- 0 --> תכנון -- 0 --> ב -- 1 --> לידה -- 0 --> מקצועית -- 0 --> ARCHITECTUR"

$ws.Range("A27").Value = "my_files\CVשירה ג'יקובס  .pdf"
$ws.Range("B27").Value = "בתקשורת -- 0 --> תכנון -- 0 --> ב -- 1 --> לידה -- 0 --> מקצועית -- 0 --> ARCHITECTUR"

$ws.Range("A28").Value = "my_files\קורות חיים אורי עוז מרזם.pdf"
$ws.Range("B28").Value = "בתקשורת -- 0 --> תכנון -- 0 --> ב -- 1 --> לידה -- 0 --> מקצועית -- 0 --> ARCHITECTUR"
